$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(0.040213000367242, 0.0194638266617701, 0.00367242012486228, 0.999081894968784, 0, 0.000734484024972457, 0.000918105031215571, 0, 0.998898273962541, 0.000183621006243114, 0.000550863018729343, 0.999081894968784, 0.997980168931326, 0.000734484024972457, 0.999449136981271, 0.998714652956298, 0.00110172603745869, 0.0328681601175174, 0.00146896804994491, 0.0146896804994491, 0.986044803525523, 0.00697759823723834, 0.972824091076019)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, 2+$i).Value = $row2[$i]
}

$row3 = @(0.936834373852369, 0.976312890194638, 0.0609621740727139, 0, 0.000367242012486228, 0, 0, 0.999816378993757, 0, 0, 0, 0.000367242012486228, 0, 0, 0, 0, 0, 0.937752478883584, 0, 0, 0.000367242012486228, 0, 0.00477414616232097)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, 2+$i).Value = $row3[$i]
}

$row4 = @(0.00110172603745869, 0.00183621006243114, 0, 0.000550863018729343, 0.000550863018729343, 0.998898273962541, 0.998714652956298, 0, 0.00110172603745869, 0.999632757987514, 0.999081894968784, 0.000550863018729343, 0.00201983106867426, 0.999081894968784, 0.000550863018729343, 0.00110172603745869, 0.998347410943812, 0.000183621006243114, 0.997980168931326, 0.984392214469335, 0.012853470437018, 0.992471538744032, 0.0119353654058024)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, 2+$i).Value = $row4[$i]
}

$row5 = @(0.0211164157179581, 0.0012853470437018, 0.933896437752479, 0.000183621006243114, 0.998714652956298, 0.000367242012486228, 0.000367242012486228, 0.000183621006243114, 0, 0, 0, 0, 0, 0, 0, 0.000183621006243114, 0.000183621006243114, 0.027359529930224, 0.000183621006243114, 0.000183621006243114, 0.000367242012486228, 0, 0.00863018729342637)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, 2+$i).Value = $row5[$i]
}
